$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry describes a cell whose text content changed in the latest
# cryptos-list refresh. Values are written as literal text (matching the
# original inline-string cells) rather than being auto-converted to
# numbers by Excel's smart-entry parsing.
$updates = @(
    @{ Addr = 'D2'; Value = '27.756.51' },
    @{ Addr = 'E2'; Value = '  +6.33%  ' },
    @{ Addr = 'D3'; Value = '1.737.57' },
    @{ Addr = 'E3'; Value = '  +5.11%  ' },
    @{ Addr = 'D4'; Value = '1.004' },
    @{ Addr = 'E4'; Value = '  -0.04%  ' },
    @{ Addr = 'D5'; Value = '227.69' },
    @{ Addr = 'E5'; Value = '  +4.06%  ' },
    @{ Addr = 'D6'; Value = '0.5461' },
    @{ Addr = 'E6'; Value = '  +3.99%  ' },
    @{ Addr = 'D7'; Value = '1.004' },
    @{ Addr = 'D8'; Value = '0.2761' },
    @{ Addr = 'E8'; Value = '  +3.17%  ' },
    @{ Addr = 'D9'; Value = '0.06714' },
    @{ Addr = 'E9'; Value = '  +5.52%  ' },
    @{ Addr = 'D10'; Value = '21.93' },
    @{ Addr = 'E10'; Value = '  +6.82%  ' },
    @{ Addr = 'D11'; Value = '0.07785' },
    @{ Addr = 'E11'; Value = '  +1.16%  ' },
    @{ Addr = 'D12'; Value = '4.692' },
    @{ Addr = 'B13'; Value = 'WrappedEther' },
    @{ Addr = 'C13'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' },
    @{ Addr = 'D13'; Value = '1.753.31' },
    @{ Addr = 'E13'; Value = '  +5.85%  ' },
    @{ Addr = 'B14'; Value = 'WrappedliquidstakedEther2.0' },
    @{ Addr = 'C14'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth' },
    @{ Addr = 'D14'; Value = '1.976.14' },
    @{ Addr = 'E14'; Value = '  +5.08%  ' },
    @{ Addr = 'E15'; Value = '  +6.60%  ' },
    @{ Addr = 'D16'; Value = '0.0₅8431' },
    @{ Addr = 'E16'; Value = '  +2.50%  ' },
    @{ Addr = 'D17'; Value = '69.48' },
    @{ Addr = 'E17'; Value = '  +5.99%  ' },
    @{ Addr = 'D18'; Value = '27.768.19' },
    @{ Addr = 'E18'; Value = '  +6.37%  ' },
    @{ Addr = 'D19'; Value = '226.96' },
    @{ Addr = 'E19'; Value = '  +18.71%  ' },
    @{ Addr = 'D20'; Value = '4.840' },
    @{ Addr = 'E20'; Value = '  +3.26%  ' },
    @{ Addr = 'D21'; Value = '1.003' },
    @{ Addr = 'E21'; Value = '  -0.07%  ' },
    @{ Addr = 'D22'; Value = '10.90' },
    @{ Addr = 'E22'; Value = '  +5.48%  ' },
    @{ Addr = 'D23'; Value = '6.232' },
    @{ Addr = 'E23'; Value = '  +4.33%  ' },
    @{ Addr = 'E24'; Value = '  -0.15%  ' },
    @{ Addr = 'D25'; Value = '147.12' },
    @{ Addr = 'E25'; Value = '  +0.72%  ' },
    @{ Addr = 'D26'; Value = '0.1252' },
    @{ Addr = 'E26'; Value = '  +4.08%  ' },
    @{ Addr = 'D27'; Value = '1.707' },
    @{ Addr = 'E27'; Value = '  +12.46%  ' },
    @{ Addr = 'B28'; Value = 'EthereumClassic' },
    @{ Addr = 'C28'; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc' },
    @{ Addr = 'D28'; Value = '17.20' },
    @{ Addr = 'E28'; Value = '  +7.69%  ' },
    @{ Addr = 'B29'; Value = 'Cosmos' },
    @{ Addr = 'C29'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' },
    @{ Addr = 'D29'; Value = '7.458' },
    @{ Addr = 'E29'; Value = '  +2.80%  ' },
    @{ Addr = 'D30'; Value = '0.05680' },
    @{ Addr = 'E30'; Value = '  +0.60%  ' },
    @{ Addr = 'D31'; Value = '1.314' },
    @{ Addr = 'E31'; Value = '  +3.30%  ' },
    @{ Addr = 'D32'; Value = '3.701' },
    @{ Addr = 'E32'; Value = '  +5.93%  ' },
    @{ Addr = 'D33'; Value = '3.518' },
    @{ Addr = 'E33'; Value = '  +4.00%  ' },
    @{ Addr = 'D34'; Value = '1.686' },
    @{ Addr = 'E34'; Value = '  +6.74%  ' },
    @{ Addr = 'D35'; Value = '0.9763' },
    @{ Addr = 'E35'; Value = '  +3.19%  ' },
    @{ Addr = 'D36'; Value = '2.857' },
    @{ Addr = 'E36'; Value = '  +2.15%  ' },
    @{ Addr = 'D37'; Value = '2.449' },
    @{ Addr = 'E37'; Value = '  +1.65%  ' },
    @{ Addr = 'D38'; Value = '0.5974' },
    @{ Addr = 'E38'; Value = '  +3.31%  ' },
    @{ Addr = 'D40'; Value = '5.904' },
    @{ Addr = 'E40'; Value = '  -1.16%  ' },
    @{ Addr = 'B41'; Value = 'Maker' },
    @{ Addr = 'C41'; Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr' },
    @{ Addr = 'D41'; Value = '1.051.16' },
    @{ Addr = 'E41'; Value = '  +2.97%  ' },
    @{ Addr = 'B42'; Value = 'TrustWalletToken' },
    @{ Addr = 'C42'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' },
    @{ Addr = 'D42'; Value = '0.8479' },
    @{ Addr = 'E42'; Value = '  +0.29%  ' },
    @{ Addr = 'D43'; Value = '1.004' },
    @{ Addr = 'E43'; Value = '  -0.03%  ' },
    @{ Addr = 'D44'; Value = '102.17' },
    @{ Addr = 'E44'; Value = '  +0.88%  ' },
    @{ Addr = 'D45'; Value = '1.880.56' },
    @{ Addr = 'E45'; Value = '  +4.98%  ' },
    @{ Addr = 'D46'; Value = '0.0₈117' },
    @{ Addr = 'E46'; Value = '  +14.06%  ' },
    @{ Addr = 'D47'; Value = '59.53' },
    @{ Addr = 'E47'; Value = '  +1.94%  ' },
    @{ Addr = 'D48'; Value = '8.280' },
    @{ Addr = 'E48'; Value = '  +2.79%  ' },
    @{ Addr = 'D49'; Value = '0.4438' },
    @{ Addr = 'E49'; Value = '  +2.20%  ' },
    @{ Addr = 'D50'; Value = '1.002' },
    @{ Addr = 'E50'; Value = '  -0.09%  ' },
    @{ Addr = 'D51'; Value = '0.05316' },
    @{ Addr = 'E51'; Value = '  -0.49%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.ClearFormats()
}
